$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-DashboardRow {
    param($Row, $NewDate, $NewQ)

    # Read current Q,R,S,T values before overwriting, so we can shift them right
    # into R,S,T,U (the oldest value in U is dropped).
    $oldQ = $ws.Range("Q$Row").Value2
    $oldR = $ws.Range("R$Row").Value2
    $oldS = $ws.Range("S$Row").Value2
    $oldT = $ws.Range("T$Row").Value2

    $ws.Range("N$Row").Value = $NewDate

    $ws.Range("Q$Row").Value = $NewQ
    $ws.Range("R$Row").Value = $oldQ
    $ws.Range("S$Row").Value = $oldR
    $ws.Range("T$Row").Value = $oldS
    $ws.Range("U$Row").Value = $oldT
}

# Row 29: T5YIFR
Update-DashboardRow 29 46024 2.22

# Row 30: T10YIE
Update-DashboardRow 30 46024 2.25

# Row 47: DFF - only the date advances, the daily values stay the same
$ws.Range("N47").Value = 46023

# Row 48: DGS2
Update-DashboardRow 48 46022 3.47

# Row 49: DGS5
Update-DashboardRow 49 46022 3.73

# Row 50: DGS10
Update-DashboardRow 50 46022 4.18
